# Add a new "MANPOWER UK HOLDINGS LIMITED" lookup result as row 10 of the
# VAT-number / company-name matching report, resize the Company Name /
# ReturnedCompany columns so the longer values are fully visible, and move
# the active selection down to below the newly-added row (matches the
# "Updated requirements with output example and name matching explanation"
# commit).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New example row: a VAT number that is already present in the sheet
# (492 4357 26) looked up against "MANPOWER UK HOLDINGS LIMITED" - the VAT
# lookup succeeds (Valid = TRUE) but returns a different registered company
# name (RANDSTAD UK HOLDING LIMITED), so NamesMatch = FALSE.
$ws.Range("A10").Value = "492 4357 26"
$ws.Range("B10").Value = "MANPOWER UK HOLDINGS LIMITED"
$ws.Range("C10").Value = $true
$ws.Range("D10").Value = "RANDSTAD UK HOLDING LIMITED"
$ws.Range("E10").Value = $false
$ws.Range("F10").Value = 75

# Colour the new row the same "red" used elsewhere in the sheet for rows
# where the returned company name does not match the searched name.
$ws.Range("A10:F10").Interior.Color = 6053069

# Company Name (B) and ReturnedCompany (D) now hold longer values - widen
# both columns so the text fits (bestFit), like the rest of the table.
$ws.Columns.Item(2).AutoFit()
$ws.Columns.Item(4).AutoFit()

# Move the selection past the table, under the new last row.
[void]$ws.Range("G16").Select()

"done"
